# Generate Report for Handoff
# - Remove the row describing c5b0ffec-4d99-4a83-8572-8a84215fda90.md from every sheet
# - Update the 8b862783-29c8-4aee-95e8-a88cb7c712bf.md row to reflect a fresh handoff
#   ("Ready for handoff") with updated timestamps, and flag the zh-cn handback as stale.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$notLatestMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3cb45ca8f16f9aa1daafc57c48f3e3eed30ffbac/e2e/8b862783-29c8-4aee-95e8-a88cb7c712bf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64aa463e380bba7d808a58e2f1fc014e8c7d58d1/e2e/8b862783-29c8-4aee-95e8-a88cb7c712bf.md."

# ---- Update the 8b862783 row (row 3) on each sheet ----

# Overview: columns E (zh-cn status), F (de-de status), G (Latest HO Xliff Generate Date)
$overview.Cells.Item(3, 5).Value = "Ready for handoff"
$overview.Cells.Item(3, 6).Value = "Ready for handoff"
$overview.Cells.Item(3, 7).Value = "2016-10-24 08:05:18"

# zh-cn: column C (Status), column H (Latest Handoff Datetime), column P (Error Detail)
$zhcn.Cells.Item(3, 3).Value = "Ready for handoff"
$zhcn.Cells.Item(3, 8).Value = "2016-10-24 08:05:05"
$zhcn.Cells.Item(3, 16).Value = $notLatestMsg

# de-de: column C (Status), column H (Latest Handoff Datetime)
$dede.Cells.Item(3, 3).Value = "Ready for handoff"
$dede.Cells.Item(3, 8).Value = "2016-10-24 08:05:18"

# ---- Remove the c5b0ffec row (row 4) from each sheet, rebuilding hyperlinks ----

function Remove-HandbackRow($ws, $lastCol) {
    $kept = @()
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Row
        if ($addr -ne 4) {
            $kept += , @($hl.Range.Address(), $hl.Address, $hl.TextToDisplay)
        }
    }
    $ws.Hyperlinks.Delete()
    $ws.Rows.Item(4).Delete()
    foreach ($item in $kept) {
        $ws.Hyperlinks.Add($ws.Range($item[0]), $item[1], "", "", $item[2])
    }
}

Remove-HandbackRow $overview 7
Remove-HandbackRow $zhcn 16
Remove-HandbackRow $dede 16

# ---- Column P width tweak (zh-cn / de-de) to fit the longer Error Detail text ----
$zhcn.Columns.Item(16).ColumnWidth = 40
$dede.Columns.Item(16).ColumnWidth = 40
